$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 7).Value = 1.48
$ws.Cells.Item(2, 9).Value = 6.25
$ws.Cells.Item(2, 10).Value = 2
$ws.Cells.Item(2, 12).Value = 6
$ws.Cells.Item(2, 13).Value = 1.04
$ws.Cells.Item(2, 15).Value = 1.22
$ws.Cells.Item(2, 19).Value = 2.75
$ws.Cells.Item(2, 20).Value = 1.44
$ws.Cells.Item(2, 23).Value = 1.8
$ws.Cells.Item(2, 24).Value = 1.95
$ws.Cells.Item(2, 32).Value = 8.5
$ws.Cells.Item(2, 33).Value = 17

# Row 3
$ws.Cells.Item(3, 13).Value = 1.11
$ws.Cells.Item(3, 15).Value = 1.5
$ws.Cells.Item(3, 16).Value = 2.63
$ws.Cells.Item(3, 20).Value = 1.17

# Row 4
$ws.Cells.Item(4, 13).Value = 1.2
$ws.Cells.Item(4, 14).Value = 4.33
$ws.Cells.Item(4, 15).Value = 1.8
$ws.Cells.Item(4, 16).Value = 1.91
$ws.Cells.Item(4, 21).Value = 1.8
$ws.Cells.Item(4, 22).Value = 2

# Row 5
$ws.Cells.Item(5, 7).Value = 2.05
$ws.Cells.Item(5, 9).Value = 3.9
$ws.Cells.Item(5, 12).Value = 4.75
$ws.Cells.Item(5, 26).Value = 8.5
$ws.Cells.Item(5, 28).Value = 19
$ws.Cells.Item(5, 36).Value = 8

# Row 7
$ws.Cells.Item(7, 7).Value = 1.91
$ws.Cells.Item(7, 8).Value = 3.1
$ws.Cells.Item(7, 9).Value = 4.15
$ws.Cells.Item(7, 10).Value = 2.6
$ws.Cells.Item(7, 11).Value = 1.91
$ws.Cells.Item(7, 12).Value = 4.7
$ws.Cells.Item(7, 17).Value = 2.42
$ws.Cells.Item(7, 19).Value = 4.2
$ws.Cells.Item(7, 21).Value = 1.53
$ws.Cells.Item(7, 22).Value = 2.2
$ws.Cells.Item(7, 25).Value = 5.1
$ws.Cells.Item(7, 26).Value = 7.5
$ws.Cells.Item(7, 27).Value = 9.25
$ws.Cells.Item(7, 28).Value = 16
$ws.Cells.Item(7, 29).Value = 19.5
$ws.Cells.Item(7, 36).Value = 8.75
$ws.Cells.Item(7, 37).Value = 21
$ws.Cells.Item(7, 38).Value = 15
$ws.Cells.Item(7, 39).Value = 70
$ws.Cells.Item(7, 40).Value = 50
$ws.Cells.Item(7, 41).Value = 70

# Row 8
$ws.Cells.Item(8, 12).Value = 6
$ws.Cells.Item(8, 25).Value = 5
$ws.Cells.Item(8, 31).Value = 6.5
$ws.Cells.Item(8, 33).Value = 23
$ws.Cells.Item(8, 37).Value = 23
$ws.Cells.Item(8, 39).Value = 51
$ws.Cells.Item(8, 41).Value = 51

# Row 9
$ws.Cells.Item(9, 7).Value = 2.88
$ws.Cells.Item(9, 13).Value = 1.06
$ws.Cells.Item(9, 15).Value = 1.36
$ws.Cells.Item(9, 20).Value = 1.22

# Row 10
$ws.Cells.Item(10, 7).Value = 1.95
$ws.Cells.Item(10, 9).Value = 4.2
$ws.Cells.Item(10, 13).Value = 1.1
$ws.Cells.Item(10, 15).Value = 1.5
$ws.Cells.Item(10, 20).Value = 1.14
$ws.Cells.Item(10, 26).Value = 8
$ws.Cells.Item(10, 29).Value = 21
$ws.Cells.Item(10, 36).Value = 8.5
$ws.Cells.Item(10, 37).Value = 19

# Row 11
$ws.Cells.Item(11, 7).Value = 3.25
$ws.Cells.Item(11, 8).Value = 3
$ws.Cells.Item(11, 9).Value = 2.35
$ws.Cells.Item(11, 11).Value = 1.91
$ws.Cells.Item(11, 12).Value = 3.25
$ws.Cells.Item(11, 14).Value = 7
$ws.Cells.Item(11, 15).Value = 1.5
$ws.Cells.Item(11, 16).Value = 2.5
$ws.Cells.Item(11, 17).Value = 2.6
$ws.Cells.Item(11, 18).Value = 1.48
$ws.Cells.Item(11, 19).Value = 5.5
$ws.Cells.Item(11, 20).Value = 1.14
$ws.Cells.Item(11, 21).Value = 1.57
$ws.Cells.Item(11, 22).Value = 2.25
$ws.Cells.Item(11, 23).Value = 2.1
$ws.Cells.Item(11, 24).Value = 1.67
$ws.Cells.Item(11, 25).Value = 7.5
$ws.Cells.Item(11, 26).Value = 15
$ws.Cells.Item(11, 30).Value = 41
$ws.Cells.Item(11, 31).Value = 6.5
$ws.Cells.Item(11, 33).Value = 19
$ws.Cells.Item(11, 34).Value = 67
$ws.Cells.Item(11, 36).Value = 6
$ws.Cells.Item(11, 38).Value = 10
$ws.Cells.Item(11, 39).Value = 23
$ws.Cells.Item(11, 40).Value = 23
$ws.Cells.Item(11, 42).Value = 2
$ws.Cells.Item(11, 43).Value = 1.85
$ws.Cells.Item(11, 44).Value = 4.4
$ws.Cells.Item(11, 45).Value = 1.21

# Row 12
$ws.Cells.Item(12, 17).Value = 1.63

# Row 15
$ws.Cells.Item(15, 8).Value = 4.75
$ws.Cells.Item(15, 9).Value = 7
$ws.Cells.Item(15, 12).Value = 7
$ws.Cells.Item(15, 14).Value = 17
$ws.Cells.Item(15, 17).Value = 1.6
$ws.Cells.Item(15, 18).Value = 2.3
$ws.Cells.Item(15, 32).Value = 9.5
$ws.Cells.Item(15, 41).Value = 51
$ws.Cells.Item(15, 44).Value = 2
$ws.Cells.Item(15, 45).Value = 1.85

# Row 17
$ws.Cells.Item(17, 13).Value = 1.02
$ws.Cells.Item(17, 14).Value = 21
$ws.Cells.Item(17, 15).Value = 1.13
$ws.Cells.Item(17, 16).Value = 6
$ws.Cells.Item(17, 17).Value = 1.44
$ws.Cells.Item(17, 18).Value = 2.7
$ws.Cells.Item(17, 19).Value = 2
$ws.Cells.Item(17, 20).Value = 1.73
$ws.Cells.Item(17, 25).Value = 19
$ws.Cells.Item(17, 27).Value = 15

# Row 18
$ws.Cells.Item(18, 7).Value = 5.75
$ws.Cells.Item(18, 10).Value = 5.5
$ws.Cells.Item(18, 11).Value = 2.75
$ws.Cells.Item(18, 12).Value = 1.91
$ws.Cells.Item(18, 14).Value = 23
$ws.Cells.Item(18, 15).Value = 1.1
$ws.Cells.Item(18, 16).Value = 7
$ws.Cells.Item(18, 17).Value = 1.36
$ws.Cells.Item(18, 18).Value = 3.1
$ws.Cells.Item(18, 19).Value = 1.83
$ws.Cells.Item(18, 20).Value = 1.83
$ws.Cells.Item(18, 21).Value = 1.22
$ws.Cells.Item(18, 22).Value = 4
$ws.Cells.Item(18, 23).Value = 1.5
$ws.Cells.Item(18, 24).Value = 2.5
$ws.Cells.Item(18, 25).Value = 23
$ws.Cells.Item(18, 28).Value = 67
$ws.Cells.Item(18, 29).Value = 41
$ws.Cells.Item(18, 30).Value = 34
$ws.Cells.Item(18, 31).Value = 23
$ws.Cells.Item(18, 34).Value = 41
$ws.Cells.Item(18, 35).Value = 101
$ws.Cells.Item(18, 36).Value = 12
$ws.Cells.Item(18, 37).Value = 10
$ws.Cells.Item(18, 39).Value = 12

# Row 19
$ws.Cells.Item(19, 7).Value = 2.63
$ws.Cells.Item(19, 9).Value = 2.45
$ws.Cells.Item(19, 12).Value = 3
$ws.Cells.Item(19, 25).Value = 13
$ws.Cells.Item(19, 39).Value = 26

# Row 20
$ws.Cells.Item(20, 7).Value = 2.55
$ws.Cells.Item(20, 9).Value = 2.8
$ws.Cells.Item(20, 10).Value = 3.1
$ws.Cells.Item(20, 12).Value = 3.4
$ws.Cells.Item(20, 13).Value = 1.06
$ws.Cells.Item(20, 14).Value = 10
$ws.Cells.Item(20, 25).Value = 9
$ws.Cells.Item(20, 28).Value = 23
$ws.Cells.Item(20, 33).Value = 13
$ws.Cells.Item(20, 36).Value = 10
$ws.Cells.Item(20, 37).Value = 15
$ws.Cells.Item(20, 39).Value = 29

# Row 21
$ws.Cells.Item(21, 7).Value = 3
$ws.Cells.Item(21, 8).Value = 3
$ws.Cells.Item(21, 9).Value = 2.55
$ws.Cells.Item(21, 10).Value = 3.5
$ws.Cells.Item(21, 12).Value = 3.1
$ws.Cells.Item(21, 17).Value = 2.03
$ws.Cells.Item(21, 18).Value = 1.83
$ws.Cells.Item(21, 19).Value = 3.4
$ws.Cells.Item(21, 20).Value = 1.3
$ws.Cells.Item(21, 21).Value = 1.4
$ws.Cells.Item(21, 22).Value = 2.75
$ws.Cells.Item(21, 23).Value = 1.7
$ws.Cells.Item(21, 24).Value = 2.05
$ws.Cells.Item(21, 25).Value = 10
$ws.Cells.Item(21, 30).Value = 29
$ws.Cells.Item(21, 31).Value = 9.5
$ws.Cells.Item(21, 39).Value = 23

# Row 22
$ws.Cells.Item(22, 9).Value = 3.7
$ws.Cells.Item(22, 26).Value = 9.5
$ws.Cells.Item(22, 28).Value = 19

# Row 23
$ws.Cells.Item(23, 15).Value = 1.36
$ws.Cells.Item(23, 16).Value = 3
$ws.Cells.Item(23, 17).Value = 2.2
$ws.Cells.Item(23, 18).Value = 1.65

# Row 24
$ws.Cells.Item(24, 7).Value = 1.67
$ws.Cells.Item(24, 8).Value = 3.7
$ws.Cells.Item(24, 9).Value = 5.25
$ws.Cells.Item(24, 10).Value = 2.38
$ws.Cells.Item(24, 25).Value = 5
$ws.Cells.Item(24, 26).Value = 6.5
$ws.Cells.Item(24, 40).Value = 51

# Row 26
$ws.Cells.Item(26, 7).Value = 3.3
$ws.Cells.Item(26, 10).Value = 4
$ws.Cells.Item(26, 11).Value = 1.95
$ws.Cells.Item(26, 13).Value = 1.1
$ws.Cells.Item(26, 14).Value = 7
$ws.Cells.Item(26, 15).Value = 1.4
$ws.Cells.Item(26, 16).Value = 2.75
$ws.Cells.Item(26, 17).Value = 2.35
$ws.Cells.Item(26, 18).Value = 1.57
$ws.Cells.Item(26, 19).Value = 4.33
$ws.Cells.Item(26, 20).Value = 1.2
$ws.Cells.Item(26, 21).Value = 1.53
$ws.Cells.Item(26, 22).Value = 2.38
$ws.Cells.Item(26, 23).Value = 2
$ws.Cells.Item(26, 24).Value = 1.73
$ws.Cells.Item(26, 25).Value = 8.5
$ws.Cells.Item(26, 27).Value = 13
$ws.Cells.Item(26, 31).Value = 7
$ws.Cells.Item(26, 33).Value = 17
$ws.Cells.Item(26, 35).Value = 351
$ws.Cells.Item(26, 36).Value = 6.5
$ws.Cells.Item(26, 38).Value = 10
$ws.Cells.Item(26, 42).Value = 1.78
$ws.Cells.Item(26, 43).Value = 2.1

# Row 27
$ws.Cells.Item(27, 7).Value = 1.95
$ws.Cells.Item(27, 10).Value = 2.6
$ws.Cells.Item(27, 11).Value = 2.3
$ws.Cells.Item(27, 13).Value = 1.04
$ws.Cells.Item(27, 14).Value = 13
$ws.Cells.Item(27, 15).Value = 1.2
$ws.Cells.Item(27, 16).Value = 4.33
$ws.Cells.Item(27, 17).Value = 1.67
$ws.Cells.Item(27, 18).Value = 2.15
$ws.Cells.Item(27, 19).Value = 2.63
$ws.Cells.Item(27, 20).Value = 1.44
$ws.Cells.Item(27, 21).Value = 1.33
$ws.Cells.Item(27, 22).Value = 3.25
$ws.Cells.Item(27, 25).Value = 9.5
$ws.Cells.Item(27, 27).Value = 9
$ws.Cells.Item(27, 29).Value = 15
$ws.Cells.Item(27, 31).Value = 13
$ws.Cells.Item(27, 32).Value = 7
$ws.Cells.Item(27, 35).Value = 151
$ws.Cells.Item(27, 36).Value = 13
$ws.Cells.Item(27, 44).Value = 1.97
$ws.Cells.Item(27, 45).Value = 1.77
